$wb = $excel.ActiveWorkbook

# New result values (row 2, columns A-O) for each sheet, as scraped from server.
# Sheets are processed in workbook order: 2025, 2030, 2035, 2040, 2045, 2050.
$sheetRows = @(
    @(0, 240.6633596400077, 0, 0, 29093.48942174729, 0, 8095.925712661834, 0, 14865.25872276, 0, 0, 50912.59821312752, 11247.09127927, 7296.869202533937, 6890.515200515633),
    @(0, 5552.285197058362, 0, 0, 56083.30735895225, 0, 8095.925712661834, 0, 37116.14035239558, 0, 0, 73331.63626107729, 21807.55922437125, 10850.09268200976, 9413.745657571631),
    @(1201.877283394627, 7553.073389192201, 0, 0, 69964.32397151143, 0, 8095.925712661834, 0, 53248.72788943995, 0, 0, 73331.63626107729, 28263.25588588425, 15462.19557169876, 15293.66517109082),
    @(1228.359695667991, 7553.073389192201, 0, 0, 69964.32397151143, 0, 8095.925712661834, 0, 53248.72788943995, 0, 0, 73331.63626107729, 28263.25588588425, 15462.19557169876, 15293.66517109082),
    @(5377.600707020665, 7553.073389192201, 0, 0, 69964.32397151143, 0, 8095.925712661834, 0, 53248.72788943995, 0, 0, 73331.63626107729, 28263.25588588425, 15462.19557169876, 15293.66517109082),
    @(5377.600707020665, 7553.073389192201, 0, 0, 69964.32397151143, 0, 8095.925712661834, 0, 53248.72788943995, 0, 0, 73331.63626107729, 28263.25588588425, 15462.19557169876, 15293.66517109082)
)

for ($i = 0; $i -lt $sheetRows.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $rowValues = $sheetRows[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item(2, $col + 1).Value = $rowValues[$col]
    }
}
